$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "43.189.31"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.54%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.271.25"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "111.25"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "264.19"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.36%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.620"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.16%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.605"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.29%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "47.49"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0929"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.78"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("E13").Value = "  +0.87%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "15.42"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.45%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.613.95"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.06%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.852"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.52%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.283.99"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.30%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "43.105.04"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("E19").Value = "  -2.68%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.80"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.32%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "71.07"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.24%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.54"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.54%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "231.21"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.46%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.64"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.82%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("E28").Value = "  -0.93%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "40.21"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.84%  "
$ws.Range("E30").Value = "  -1.64%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "171.44"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.71%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "21.29"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.25%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0900"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.92%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.68"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.46%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.127"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.31%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.63"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0350"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  -6.66%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.78"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.93%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.05%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "76.67"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +9.66%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "13.89"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +9.79%  "
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("E46").Value = "  +0.17%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.36"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0991"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "101.14"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.48%  "
